$wb = $excel.ActiveWorkbook

# --- 1. "HSL Treinseries": sort data A2:G45 by column A ascending ---
$wsTrein = $wb.Worksheets.Item("HSL Treinseries")
$rng = $wsTrein.Range("A1:G45")
$key = $wsTrein.Range("A1:A45")
$rng.Sort($key, 1, $null, $null, 1, $null, $null, 1)

# --- 2. Add ratio "0.05" next to the Eurostar (9100) row, column G ---
$wsTrein.Range("G40").NumberFormat = "@"
$wsTrein.Range("G40").Value = "0.05"

# --- 3. Update selection on "HSL Treinseries" ---
$wsTrein.Range("E9").Select()

# --- 4. "HSL Scenarios": expand the AI description text for scenario 5801 ---
$wsScen = $wb.Worksheets.Item("HSL Scenarios")
$wsScen.Range("C2").Value = "Hfd:`n- 4640`n- 73140 (Hfdo212)`nHfd-Shl:`n- 4340 (Hfdo213)`n- 9100/9300 niet ivm tunnel`n- 3240 niet ivm tunnel`n- 3341`nShl:`n- 1040"

# --- 5. Make "HSL Scenarios" the active sheet/tab ---
$wsScen.Select()
